$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duration_Seconds (column C) values for rows 3,5,7,...,41 change from 5 to 2
$rows = @(3,5,7,9,11,13,15,17,19,21,23,25,27,29,31,33,35,37,39,41)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 2
}

# Row 42 column C changes from 2 to 16
$ws.Cells.Item(42, 3).Value = 16

# Update the active selection to G12 (matches the saved selection in the diff)
$ws.Range("G12").Select()
